# Update the "Förändrad" (Changed) date column (C) for rows 2-66 from
# 45190 (2023-09-21) to 45192 (2023-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 66; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45190) {
        $cell.Value2 = 45192
    }
}
